$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 308.625
$ws.Range("I9").Value = 117.7
$ws.Range("J9").Value = 626.8333
$ws.Range("K9").Value = 117.7
$ws.Range("L9").Value = 626.8333
$ws.Range("M9").Value = 51.3
$ws.Range("N9").Value = -964.8333

$ws.Range("H15").Value = 2847.7632
$ws.Range("I15").Value = 2847.7632
$ws.Range("K15").Value = 8543.2896
$ws.Range("M15").Value = -8374.2896

$ws.Range("H43").Value = 7856.375
$ws.Range("I43").Value = 6790
$ws.Range("K43").Value = 6790
$ws.Range("M43").Value = -6721

$ws.Range("H64").Value = 100007040
$ws.Range("I64").Value = 166671730
$ws.Range("K64").Value = 166671730
$ws.Range("M64").Value = -166671482

$ws.Range("H67").Value = 100007040
$ws.Range("I67").Value = 166671730
$ws.Range("K67").Value = 166671730
$ws.Range("M67").Value = -166670872

$ws.Range("H74").Value = 5531.3335
$ws.Range("I74").Value = 4547
$ws.Range("K74").Value = 4547
$ws.Range("M74").Value = -3611

$ws.Range("H76").Value = 6882.1
$ws.Range("I76").Value = 4783.3335
$ws.Range("J76").Value = 10030.25
$ws.Range("K76").Value = 4783.3335
$ws.Range("L76").Value = 10030.25
$ws.Range("M76").Value = -4468.3335
$ws.Range("N76").Value = -10660.25

$ws.Range("H77").Value = 5531.3335
$ws.Range("I77").Value = 4547
$ws.Range("K77").Value = 22735
$ws.Range("M77").Value = -18055

$ws.Range("H79").Value = 6882.1
$ws.Range("I79").Value = 4783.3335
$ws.Range("J79").Value = 10030.25
$ws.Range("K79").Value = 4783.3335
$ws.Range("L79").Value = 10030.25
$ws.Range("M79").Value = -3691.3335
$ws.Range("N79").Value = -12214.25

$ws.Range("H113").Value = 11836.833
$ws.Range("I113").Value = 17275.334
$ws.Range("J113").Value = 6398.3335
$ws.Range("K113").Value = 17275.334
$ws.Range("L113").Value = 6398.3335
$ws.Range("M113").Value = -14021.334
$ws.Range("N113").Value = -12906.3335

$ws.Range("H115").Value = 368
$ws.Range("I115").Value = 381.55554
$ws.Range("J115").Value = 337.5
$ws.Range("K115").Value = 1144.66662
$ws.Range("L115").Value = 1012.5
$ws.Range("M115").Value = 422.33338
$ws.Range("N115").Value = -4146.5

$ws.Range("H127").Value = 557193.9
$ws.Range("I127").Value = 557193.9
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1671581.7
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -1666621.7
$ws.Range("N127").ClearContents()

$ws.Range("H129").Value = 1481.5454
$ws.Range("I129").Value = 897.5714
$ws.Range("K129").Value = 2692.7142
$ws.Range("M129").Value = 2307.2858

$ws.Range("H138").Value = 2738.9673
$ws.Range("I138").Value = 4650.5
$ws.Range("J138").Value = 2530.4363
$ws.Range("K138").Value = 13951.5
$ws.Range("L138").Value = 7591.3089
$ws.Range("M138").Value = -8811.5
$ws.Range("N138").Value = -17871.3089

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3180.9
$ws.Range("I32").Value = 2251.1707
$ws.Range("K32").Value = 2251.1707
$ws.Range("M32").Value = -1964.1707

$ws.Range("H102").Value = 2557
$ws.Range("I102").Value = 1636.5714
$ws.Range("K102").Value = 1636.5714
$ws.Range("M102").Value = -14.57140000000004

$ws.Range("H122").Value = 1283038.1
$ws.Range("I122").Value = 2264296.5
$ws.Range("J122").Value = 7402.2
$ws.Range("K122").Value = 6792889.5
$ws.Range("L122").Value = 22206.6
$ws.Range("M122").Value = -6790439.5
$ws.Range("N122").Value = -27106.6

$ws.Range("H133").Value = 49948.5
$ws.Range("J133").Value = 49948.5
$ws.Range("L133").Value = 49948.5
$ws.Range("N133").Value = -55008.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 21440.666
$ws.Range("I99").Value = 19893.54
$ws.Range("K99").Value = 19893.54
$ws.Range("M99").Value = -18395.54

$ws.Range("H134").Value = 37655.4
$ws.Range("I134").Value = 33876.09
$ws.Range("K134").Value = 101628.27
$ws.Range("M134").Value = -99093.26999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 14333.333
$ws.Range("J80").Value = 14333.333
$ws.Range("L80").Value = 14333.333
$ws.Range("N80").Value = -16579.333

$ws.Range("H83").Value = 14333.333
$ws.Range("J83").Value = 14333.333
$ws.Range("L83").Value = 42999.999
$ws.Range("N83").Value = -54231.999

$ws.Range("H122").Value = 2484.842
$ws.Range("I122").Value = 2072.875
$ws.Range("J122").Value = 4682
$ws.Range("K122").Value = 6218.625
$ws.Range("L122").Value = 14046
$ws.Range("M122").Value = -3768.625
$ws.Range("N122").Value = -18946

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H93").Value = 3300
$ws.Range("J93").Value = 4600
$ws.Range("L93").Value = 13800
$ws.Range("N93").Value = -17544

$ws.Range("H122").Value = 9359769
$ws.Range("J122").Value = 2028440.9
$ws.Range("L122").Value = 18255968.1
$ws.Range("N122").Value = -18260868.1

$ws.Range("H131").Value = 1462.82
$ws.Range("I131").Value = 1076.6666
$ws.Range("J131").Value = 1474.763
$ws.Range("K131").Value = 3229.9998
$ws.Range("L131").Value = 4424.289
$ws.Range("M131").Value = 1810.0002
$ws.Range("N131").Value = -14504.289

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 24201
$ws.Range("I70").Value = 23051.5
$ws.Range("K70").Value = 23051.5
$ws.Range("M70").Value = -22781.5

$ws.Range("H73").Value = 24201
$ws.Range("I73").Value = 23051.5
$ws.Range("K73").Value = 23051.5
$ws.Range("M73").Value = -22115.5

$ws.Range("H80").Value = 12696.875
$ws.Range("I80").Value = 11113.3125
$ws.Range("J80").Value = 15864
$ws.Range("K80").Value = 11113.3125
$ws.Range("L80").Value = 15864
$ws.Range("M80").Value = -10115.3125
$ws.Range("N80").Value = -17860

$ws.Range("H83").Value = 12696.875
$ws.Range("I83").Value = 11113.3125
$ws.Range("J83").Value = 15864
$ws.Range("K83").Value = 55566.5625
$ws.Range("L83").Value = 79320
$ws.Range("M83").Value = -50574.5625
$ws.Range("N83").Value = -89304

$ws.Range("H102").Value = 3475703.5
$ws.Range("I102").Value = 5957101.5
$ws.Range("J102").Value = 1745.9333
$ws.Range("K102").Value = 5957101.5
$ws.Range("L102").Value = 1745.9333
$ws.Range("M102").Value = -5955479.5
$ws.Range("N102").Value = -4989.9333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 24392346
$ws.Range("I22").Value = 1499.8636
$ws.Range("J22").Value = 52634376
$ws.Range("K22").Value = 1499.8636
$ws.Range("L22").Value = 52634376
$ws.Range("M22").Value = -1204.8636
$ws.Range("N22").Value = -52634966

$ws.Range("H27").Value = 24392346
$ws.Range("I27").Value = 1499.8636
$ws.Range("J27").Value = 52634376
$ws.Range("K27").Value = 1499.8636
$ws.Range("L27").Value = 52634376
$ws.Range("M27").Value = -1392.8636
$ws.Range("N27").Value = -52634590

$ws.Range("H40").Value = 6542026.5
$ws.Range("I40").Value = 6838.625
$ws.Range("J40").Value = 58823530
$ws.Range("K40").Value = 6838.625
$ws.Range("L40").Value = 58823530
$ws.Range("M40").Value = -6702.625
$ws.Range("N40").Value = -58823802

$ws.Range("H136").Value = 16164.634
$ws.Range("I136").Value = 19424.166
$ws.Range("J136").Value = 13991.611
$ws.Range("K136").Value = 58272.49800000001
$ws.Range("L136").Value = 41974.833
$ws.Range("M136").Value = -55722.49800000001
$ws.Range("N136").Value = -47074.833

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 694.1177
$ws.Range("I100").Value = 634.625
$ws.Range("K100").Value = 1269.25
$ws.Range("M100").Value = -728.25

$ws.Range("H122").Value = 414888.44
$ws.Range("J122").Value = 5990.1665
$ws.Range("L122").Value = 17970.4995
$ws.Range("N122").Value = -22870.4995
